$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.498.45"
$ws.Range("E2").Value = "  +0.97%  "

$ws.Range("D3").Value = "1.625.17"
$ws.Range("E3").Value = "  +1.42%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'212.34"
$ws.Range("E5").Value = "  +0.05%  "

$ws.Range("E6").Value = "  -0.07%  "

$ws.Range("D7").Value = "'0.486"
$ws.Range("E7").Value = "  +0.38%  "

$ws.Range("D8").Value = "'0.249"
$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("E9").Value = "  +0.66%  "

$ws.Range("D10").Value = "'18.86"
$ws.Range("E10").Value = "  +3.63%  "

$ws.Range("D11").Value = "'0.0827"
$ws.Range("E11").Value = "  +2.22%  "

$ws.Range("D12").Value = "1.851.08"
$ws.Range("E12").Value = "  +1.48%  "

$ws.Range("D13").Value = "1.630.97"
$ws.Range("E13").Value = "  +1.90%  "

$ws.Range("D14").Value = "'4.03"
$ws.Range("E14").Value = "  +0.29%  "

$ws.Range("E15").Value = "  +1.61%  "

$ws.Range("D16").Value = "26.475.57"
$ws.Range("E16").Value = "  +0.99%  "

$ws.Range("D17").Value = "'62.78"
$ws.Range("E17").Value = "  +2.38%  "

$ws.Range("E18").Value = "  +0.09%  "

$ws.Range("E19").Value = "  -0.06%  "

$ws.Range("D20").Value = "'202.95"
$ws.Range("E20").Value = "  -0.16%  "

$ws.Range("D21").Value = "'4.29"
$ws.Range("E21").Value = "  +0.02%  "

$ws.Range("D22").Value = "'9.34"
$ws.Range("E22").Value = "  +0.85%  "

$ws.Range("E23").Value = "  +0.71%  "

$ws.Range("E24").Value = "  -3.87%  "

$ws.Range("D25").Value = "'145.17"
$ws.Range("E25").Value = "  +0.40%  "

$ws.Range("E26").Value = "  +0.15%  "

$ws.Range("D27").Value = "'0.119"
$ws.Range("E27").Value = "  -2.27%  "

$ws.Range("D28").Value = "'15.32"
$ws.Range("E28").Value = "  +0.88%  "

$ws.Range("E29").Value = "  +1.04%  "

$ws.Range("D30").Value = "'0.0521"
$ws.Range("E30").Value = "  +5.97%  "

$ws.Range("E31").Value = "  +0.30%  "

$ws.Range("E32").Value = "  +1.39%  "

$ws.Range("D33").Value = "'2.93"
$ws.Range("E33").Value = "  +0.40%  "

$ws.Range("E34").Value = "  +1.38%  "

$ws.Range("E35").Value = "  -0.62%  "

$ws.Range("D36").Value = "1.152.77"
$ws.Range("E36").Value = "  +0.59%  "

$ws.Range("E37").Value = "  +0.81%  "

$ws.Range("D38").Value = "'0.804"
$ws.Range("E38").Value = "  +2.31%  "

$ws.Range("E39").Value = "  -0.07%  "

$ws.Range("E40").Value = "  -0.36%  "

$ws.Range("D41").Value = "'0.500"
$ws.Range("E41").Value = "  +0.76%  "

$ws.Range("D42").Value = "'5.40"
$ws.Range("E42").Value = "  +3.62%  "

$ws.Range("D43").Value = "'0.784"
$ws.Range("E43").Value = "  +0.61%  "

$ws.Range("D44").Value = "1.761.64"
$ws.Range("E44").Value = "  +1.28%  "

$ws.Range("D45").Value = "'92.30"
$ws.Range("E45").Value = "  +0.50%  "

$ws.Range("D46").Value = "'1.54"
$ws.Range("E46").Value = "  +1.59%  "

$ws.Range("D47").Value = "0.0₆0102"
$ws.Range("E47").Value = "  +5.05%  "

$ws.Range("D48").Value = "'54.05"

$ws.Range("D49").Value = "'0.0510"
$ws.Range("E49").Value = "  +0.69%  "

$ws.Range("E50").Value = "  +0.54%  "

$ws.Range("E51").Value = "  +0.22%  "
